# Add a new "Tester" column (F) with bonus/tester values next to the
# existing employee data, matching header styling used by the other
# header cells (bold font, no border) and leave the selection where the
# user would naturally end up after entering the last value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "Tester" - bold like the other headers, but without the
# border that the original header style (s=1) carries.
$ws.Range("F1").Value = "Tester"
$ws.Range("F1").Font.Bold = $true

# Data values for the new column, rows 2-6.
$ws.Range("F2").Value = 4000
$ws.Range("F3").Value = 3000
$ws.Range("F4").Value = 2000
$ws.Range("F5").Value = 1000
$ws.Range("F6").Value = 600

# Column E picked up an explicit width in the saved file.
$ws.Range("E1").ColumnWidth = 10.140625

# Leave the active selection on F7 (just below the data that was typed),
# matching where Excel leaves the cursor after entering a column of data.
$ws.Range("F7").Select() | Out-Null
